# Update "想去人数" (want-to-go count) values in the F column
# on the "展览" and "全部类型" worksheets.

$wb = $excel.ActiveWorkbook

$updates = @{
    "展览"     = @{ 2 = 6736; 15 = 1460; 16 = 21; 17 = 3386; 21 = 2017; 22 = 129; 23 = 1; 24 = 29 }
    "全部类型" = @{ 2 = 6736; 16 = 1460; 17 = 21; 18 = 3386; 22 = 2017; 23 = 129; 24 = 1; 25 = 29 }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $updates[$sheetName]
    foreach ($row in $rows.Keys) {
        $ws.Range("F$row").Value = $rows[$row]
    }
}
